$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Final attendance data (already in ascending EventDate order) ---
# row => Event name, EventDate (serial), Registered, Attended
$ws.Range("A2").Value = "SQL Saturday Jacksonville 2022"
$ws.Range("B2").Value = 44695
$ws.Range("C2").Value = 406
$ws.Range("D2").Value = 252

$ws.Range("A3").Value = "SQL Saturday New Jersey 2022"
$ws.Range("B3").Value = 44769
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()

$ws.Range("A4").Value = "SQL Saturday Baton Rouge 2022"
$ws.Range("B4").Value = 44779
$ws.Range("C4").Value = 483
$ws.Range("D4").Value = 183

$ws.Range("A5").Value = "SQL Saturday Los Angeles 2022"
$ws.Range("B5").Value = 44786
$ws.Range("C5").Value = 209
$ws.Range("D5").Value = 104

$ws.Range("A6").Value = "SQL Saturday Denver 2022"
$ws.Range("B6").Value = 44821
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()

$ws.Range("A7").Value = "SQL Saturday Boston 2022"
$ws.Range("B7").Value = 44842
$ws.Range("C7").Value = 283
$ws.Range("D7").Value = 142

$ws.Range("A8").Value = "SQL Saturday Orlando 2022"
$ws.Range("B8").Value = 44842
$ws.Range("C8").Value = 240
$ws.Range("D8").Value = 110

$ws.Range("A9").Value = "SQL Saturday Memphis 2022"
$ws.Range("B9").Value = 44849
$ws.Range("C9").Value = 75
$ws.Range("D9").Value = 30

$ws.Range("A10").Value = "SQL Saturday Toronto 2022"
$ws.Range("B10").Value = 44856
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = 115

$ws.Range("A11").Value = "SQL Saturday Richmond 2022"
$ws.Range("B11").Value = 44856
$ws.Range("C11").Value = 115
$ws.Range("D11").Value = 60

# Rows 12:18 stay empty in columns A:D (already blank)

# --- Recompute the "No show rate" formula for rows that have data ---
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("E$r").Formula = "=+(C$r-D$r)/C$r"
}

# --- Rows 12:18 have no event data any more, so No show rate is blank ---
for ($r = 12; $r -le 18; $r++) {
    $ws.Range("E$r").ClearContents()
}

# --- Bold the header row ---
$ws.Range("A1:E1").Font.Bold = $true

# --- Record a sort-by-date state (matches author re-sorting the table) ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B18"))
$ws.Sort.SetRange($ws.Range("A2:E18"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection / active cell ---
$ws.Range("H10").Select()

Write-Host "Edit applied"
